# feat: add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new top data row for the 2022-Q4
#    quarter, pushing the existing quarter rows down by one and
#    renumbering the leading index column.
# 2) Insert a new "2022-Q4" worksheet (cloned from the "2022-Q3" sheet so
#    it keeps the same headers/styles/column widths) right after "总计",
#    then overwrite its single data row with the new quarter's figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 sheet: insert new row 2 for "2022-Q4"
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Copy row 3's formatting (now the old row 2, "2022-Q3") onto the new
# blank row 2 so the new row matches the existing look (bold/bordered
# index style in column A, plain in B:D) instead of Excel's default
# insert formatting.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.06

# Renumber the leading index column (0,1,2,3,4,5) for the rows that got
# pushed down.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, positioned right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $summary)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fund code/name stay the same (010764 / 九泰锐升混合); only the
# quarterly figures change. D:G are stored as text in the source data,
# so force text entry with a leading apostrophe; H is a genuine number.
$q4.Cells.Item(2, 4).Value = "'2.02"
$q4.Cells.Item(2, 5).Value = "'88.37"
$q4.Cells.Item(2, 6).Value = "'3.12"
$q4.Cells.Item(2, 7).Value = "'0.0630"
$q4.Cells.Item(2, 8).Value = 9

# The leading apostrophe marks the cells as "number stored as text" with
# a quote-prefix style; strip that back off (copy the plain/default
# format from an untouched blank cell) so D2:G2 end up with the same
# styling as the rest of the sheet.
$q4.Cells.Item(1, 20).Copy()
$q4.Range("D2:G2").PasteSpecial(-4122)
$q4.Cells.Item(1, 20).Clear()

# Restore the original active tab ("2021-Q3") so the copy operation
# doesn't leave the new/duplicated sheet selected instead.
$wb.Worksheets.Item("2021-Q3").Activate()
